# Auto-update Cloudflare export
# Replace the 'settings' column (H) value "{'flatten_cname': False}" with "{}"
# for the specific DNS record rows that were updated in this export refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(4,5,15,16,17,18,28,33,39,40,41,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,85,86,87,88,89,90,92,93,94,95,96,112,113,114,115,116,117,118,129,130,143)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 8)  # column H = 8
    if ($cell.Text -eq "{'flatten_cname': False}") {
        $cell.Value = "{}"
    }
}
